$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H12").Value = 943.5
$ws.Range("I12").Value = 899.5
$ws.Range("K12").Value = 899.5
$ws.Range("M12").Value = -729.5
$ws.Range("H19").Value = 1100
$ws.Range("J19").Value = 1150
$ws.Range("L19").Value = 1150
$ws.Range("N19").Value = -1500
$ws.Range("H33").Value = 219
$ws.Range("I33").Value = 203.6
$ws.Range("K33").Value = 203.6
$ws.Range("M33").Value = 25.40000000000001
$ws.Range("H40").Value = 12223.792
$ws.Range("J40").Value = 14051.909
$ws.Range("L40").Value = 14051.909
$ws.Range("N40").Value = -14401.909
$ws.Range("H51").Value = 8167.4814
$ws.Range("I51").Value = 4992.5713
$ws.Range("K51").Value = 4992.5713
$ws.Range("M51").Value = -4508.5713
$ws.Range("H86").Value = 1321147.5
$ws.Range("I86").Value = 1909418.8
$ws.Range("K86").Value = 1909418.8
$ws.Range("M86").Value = -1908295.8
$ws.Range("H89").Value = 1321147.5
$ws.Range("I89").Value = 1909418.8
$ws.Range("K89").Value = 9547094
$ws.Range("M89").Value = -9541478
$ws.Range("H106").Value = 25643420
$ws.Range("I106").Value = 27779538
$ws.Range("K106").Value = 27779538
$ws.Range("M106").Value = -27778907
$ws.Range("H107").Value = 546.4167
$ws.Range("J107").Value = 766.5
$ws.Range("L107").Value = 766.5
$ws.Range("N107").Value = -4606.5
$ws.Range("H137").Value = 12720483
$ws.Range("I137").Value = 910694.6
$ws.Range("J137").Value = 20839712
$ws.Range("K137").Value = 2732083.8
$ws.Range("L137").Value = 62519136
$ws.Range("M137").Value = -2729533.8
$ws.Range("N137").Value = -62524236
$ws.Range("H138").Value = 2158.03
$ws.Range("I138").Value = 1190.3334
$ws.Range("K138").Value = 3571.0002
$ws.Range("M138").Value = 1568.9998
$ws.Range("H141").Value = 2326.5588
$ws.Range("I141").Value = 2326.5588
$ws.Range("K141").Value = 6979.676399999999
$ws.Range("M141").Value = -1799.676399999999

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 15420.807
$ws.Range("J32").Value = 17339.3
$ws.Range("L32").Value = 17339.3
$ws.Range("N32").Value = -17913.3
$ws.Range("H39").Value = 7985.125
$ws.Range("I39").Value = 6983
$ws.Range("K39").Value = 6983
$ws.Range("M39").Value = -6463
$ws.Range("H45").Value = 2641.8462
$ws.Range("I45").Value = 1418.625
$ws.Range("K45").Value = 1418.625
$ws.Range("M45").Value = -1041.625
$ws.Range("H97").Value = 274.41934
$ws.Range("I97").Value = 291.69565
$ws.Range("J97").Value = 224.75
$ws.Range("K97").Value = 291.69565
$ws.Range("L97").Value = 224.75
$ws.Range("M97").Value = 204.30435
$ws.Range("N97").Value = -1216.75
$ws.Range("H110").Value = 2555648
$ws.Range("I110").Value = 4084517
$ws.Range("K110").Value = 4084517
$ws.Range("M110").Value = -4082472
$ws.Range("H132").Value = 12083.566
$ws.Range("I132").Value = 15614
$ws.Range("K132").Value = 46842
$ws.Range("M132").Value = -44312

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 2126.5833
$ws.Range("I105").Value = 2169
$ws.Range("J105").Value = 1829.6666
$ws.Range("K105").Value = 2169
$ws.Range("L105").Value = 1829.6666
$ws.Range("M105").Value = -422
$ws.Range("N105").Value = -5323.6666
$ws.Range("H134").Value = 946.5484
$ws.Range("I134").Value = 938.1111
$ws.Range("K134").Value = 2814.3333
$ws.Range("M134").Value = -279.3332999999998

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6409.354
$ws.Range("I31").Value = 1714.6
$ws.Range("K31").Value = 1714.6
$ws.Range("M31").Value = -1419.6
$ws.Range("H34").Value = 6409.354
$ws.Range("I34").Value = 1714.6
$ws.Range("K34").Value = 1714.6
$ws.Range("M34").Value = -1512.6
$ws.Range("H38").Value = 2332.6667
$ws.Range("I38").Value = 2999.5
$ws.Range("J38").Value = 999
$ws.Range("K38").Value = 2999.5
$ws.Range("L38").Value = 999
$ws.Range("M38").Value = -2622.5
$ws.Range("N38").Value = -1753
$ws.Range("H46").Value = 2332.6667
$ws.Range("I46").Value = 2999.5
$ws.Range("J46").Value = 999
$ws.Range("K46").Value = 2999.5
$ws.Range("L46").Value = 999
$ws.Range("M46").Value = -2788.5
$ws.Range("N46").Value = -1421
$ws.Range("H86").Value = 10702.2
$ws.Range("I86").Value = 10814.8
$ws.Range("K86").Value = 10814.8
$ws.Range("M86").Value = -9691.799999999999
$ws.Range("H89").Value = 10702.2
$ws.Range("I89").Value = 10814.8
$ws.Range("K89").Value = 54074
$ws.Range("M89").Value = -48458
$ws.Range("H105").Value = 5683268
$ws.Range("I105").Value = 11364536
$ws.Range("J105").Value = 2000
$ws.Range("K105").Value = 11364536
$ws.Range("L105").Value = 2000
$ws.Range("M105").Value = -11362789
$ws.Range("N105").Value = -5494
$ws.Range("H132").Value = 15771.823
$ws.Range("I132").Value = 15771.823
$ws.Range("K132").Value = 47315.469
$ws.Range("M132").Value = -44785.469
$ws.Range("H141").Value = 76792.69
$ws.Range("J141").Value = 77892.47
$ws.Range("L141").Value = 77892.47
$ws.Range("N141").Value = -88252.47

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H14").Value = 220.15384
$ws.Range("I14").Value = 220.15384
$ws.Range("K14").Value = 660.4615200000001
$ws.Range("M14").Value = -487.4615200000001
$ws.Range("H68").Value = 1833.3334
$ws.Range("J68").Value = 1833.3334
$ws.Range("L68").Value = 5500.0002
$ws.Range("N68").Value = -7122.0002
$ws.Range("H71").Value = 1833.3334
$ws.Range("J71").Value = 1833.3334
$ws.Range("L71").Value = 16500.0006
$ws.Range("N71").Value = -24612.0006
$ws.Range("H92").Value = 1241.375
$ws.Range("J92").Value = 1247.6
$ws.Range("L92").Value = 3742.8
$ws.Range("N92").Value = -6238.799999999999
$ws.Range("H131").Value = 19286920
$ws.Range("J131").Value = 19286920
$ws.Range("L131").Value = 57860760
$ws.Range("N131").Value = -57870840

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 4686.65
$ws.Range("I126").Value = 2275.625
$ws.Range("J126").Value = 6294
$ws.Range("K126").Value = 6826.875
$ws.Range("L126").Value = 18882
$ws.Range("M126").Value = -4356.875
$ws.Range("N126").Value = -23822

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 5460.5
$ws.Range("I7").Value = 2739.7273
$ws.Range("J7").Value = 8181.273
$ws.Range("K7").Value = 2739.7273
$ws.Range("L7").Value = 8181.273
$ws.Range("M7").Value = -2627.7273
$ws.Range("N7").Value = -8405.273000000001
$ws.Range("H16").Value = 2956.8076
$ws.Range("I16").Value = 1637.238
$ws.Range("K16").Value = 1637.238
$ws.Range("M16").Value = -1467.238
$ws.Range("H126").Value = 5460.5
$ws.Range("I126").Value = 2739.7273
$ws.Range("J126").Value = 8181.273
$ws.Range("K126").Value = 8219.1819
$ws.Range("L126").Value = 24543.819
$ws.Range("M126").Value = -5749.1819
$ws.Range("N126").Value = -29483.819
$ws.Range("H132").Value = 4861.2173
$ws.Range("I132").Value = 4665.1763
$ws.Range("J132").Value = 5416.6665
$ws.Range("K132").Value = 13995.5289
$ws.Range("L132").Value = 16249.9995
$ws.Range("M132").Value = -11465.5289
$ws.Range("N132").Value = -21309.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H13").Value = 1799.5
$ws.Range("I13").Value = 400
$ws.Range("J13").Value = 2266
$ws.Range("K13").Value = 400
$ws.Range("L13").Value = 2266
$ws.Range("M13").Value = -260
$ws.Range("N13").Value = -2546
$ws.Range("H122").Value = 5661.5557
$ws.Range("I122").Value = 5524.0586
$ws.Range("K122").Value = 16572.1758
$ws.Range("M122").Value = -14122.1758
$ws.Range("H136").Value = 8169.882
$ws.Range("I136").Value = 2938.9688
$ws.Range("K136").Value = 8816.9064
$ws.Range("M136").Value = -6266.9064
